$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the D:E value range to retain text formatting so price/percentage
# strings (including number-like ones such as "1.00" or "0.100") are not
# auto-converted to numeric values by Excel's type inference.
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "42.545.51"
$ws.Range("E2").Value = "  -2.11%  "
$ws.Range("D3").Value = "2.346.94"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "323.23"
$ws.Range("E5").Value = "  +2.70%  "
$ws.Range("D6").Value = "99.65"
$ws.Range("E6").Value = "  -8.87%  "
$ws.Range("D7").Value = "0.634"
$ws.Range("E7").Value = "  -1.32%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").Value = "0.621"
$ws.Range("E9").Value = "  -2.42%  "
$ws.Range("D10").Value = "39.68"
$ws.Range("E10").Value = "  -8.02%  "
$ws.Range("D11").Value = "0.0919"
$ws.Range("E11").Value = "  -2.21%  "
$ws.Range("D12").Value = "8.37"
$ws.Range("E12").Value = "  -5.19%  "
$ws.Range("D13").Value = "0.996"
$ws.Range("E13").Value = "  -4.34%  "
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("D15").Value = "16.43"
$ws.Range("E15").Value = "  +0.53%  "
$ws.Range("D16").Value = "2.708.83"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").Value = "2.357.39"
$ws.Range("E17").Value = "  -2.76%  "
$ws.Range("D18").Value = "8.05"
$ws.Range("E18").Value = "  +11.22%  "
$ws.Range("D19").Value = "42.561.29"
$ws.Range("E19").Value = "  -1.99%  "
$ws.Range("E20").Value = "  -2.52%  "
$ws.Range("D21").Value = "75.75"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("E22").Value = "  +6.69%  "
$ws.Range("D23").Value = "265.89"
$ws.Range("E23").Value = "  +4.14%  "
$ws.Range("D24").Value = "2.30"
$ws.Range("E24").Value = "  -11.04%  "
$ws.Range("D25").Value = "9.99"
$ws.Range("E25").Value = "  +9.54%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").Value = "11.39"
$ws.Range("E27").Value = "  -5.44%  "
$ws.Range("D28").Value = "22.89"
$ws.Range("E28").Value = "  +2.22%  "
$ws.Range("E29").Value = "  +1.77%  "
$ws.Range("D30").Value = "175.32"
$ws.Range("E30").Value = "  +1.00%  "
$ws.Range("E31").Value = "  -2.34%  "
$ws.Range("D32").Value = "0.0895"
$ws.Range("E32").Value = "  -3.62%  "
$ws.Range("D33").Value = "35.08"
$ws.Range("E33").Value = "  -10.12%  "
$ws.Range("D34").Value = "5.99"
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("E35").Value = "  -0.40%  "
$ws.Range("D36").Value = "4.55"
$ws.Range("E36").Value = "  -8.78%  "
$ws.Range("D37").Value = "0.0356"
$ws.Range("E37").Value = "  -5.38%  "
$ws.Range("D38").Value = "2.94"
$ws.Range("E38").Value = "  +8.79%  "
$ws.Range("E39").Value = "  +0.86%  "
$ws.Range("E40").Value = "  -9.40%  "
$ws.Range("E41").Value = "  +0.75%  "
$ws.Range("D42").Value = "0.234"
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").Value = "69.63"
$ws.Range("E43").Value = "  -4.01%  "
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("D45").Value = "118.12"
$ws.Range("E45").Value = "  +6.28%  "
$ws.Range("D46").Value = "90.95"
$ws.Range("E46").Value = "  +30.04%  "
$ws.Range("E47").Value = "  -8.13%  "
$ws.Range("D48").Value = "5.47"
$ws.Range("E48").Value = "  -2.73%  "
$ws.Range("D49").Value = "9.13"
$ws.Range("E49").Value = "  -1.63%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "1.558.96"
$ws.Range("E50").Value = "  +4.11%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.100"
$ws.Range("E51").Value = "  -0.48%  "

# Restore the default (Normal) style so no stray number-format style is
# left behind on the cells.
$priceRange.Style = "Normal"

